$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new column (R) for year 2021 is being appended right after the existing
# 2020 column (Q). Copy the formatting of column Q into column R first so the
# new cells inherit the same number formats / fonts / borders as the rest of
# the data table, then overwrite the copied values with the real 2021 figures.
$ws.Range("Q4:Q14").Copy($ws.Range("R4:R14"))

$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 1
$ws.Range("R6").Value = 2.2
$ws.Range("R7").Value = 1.7
$ws.Range("R8").Value = "-"
$ws.Range("R9").Value = 0.3
$ws.Range("R10").Value = 1.1
$ws.Range("R11").Value = "-"
$ws.Range("R12").Value = 0.9
$ws.Range("R13").Value = 0.4
$ws.Range("R14").Value = 0.6

# Reflect the selection left behind in the saved workbook.
$ws.Range("S17").Select()
